$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.520.34'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.166.01'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '164.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.19%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -3.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.117'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.09%  '
$ws.Range("E10").Value = '  -0.29%  '
$ws.Range("E11").Value = '  +1.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.720.58'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("E13").Value = '  -0.89%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '64.538.49'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.28'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.86%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.158.96'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("E17").Value = '  -0.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '407.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.01%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("E25").Value = '  -2.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000102'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.83'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("E29").Value = '  +0.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.24'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.34'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.88'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.11%  '
$ws.Range("E33").Value = '  +0.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '156.46'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.95%  '
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("E36").Value = '  +1.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.687.06'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.96'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.75%  '
$ws.Range("E39").Value = '  -0.39%  '
$ws.Range("E40").Value = '  -1.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0618'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.50'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.90%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0258'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '291.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0984'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("E48").Value = '  -4.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.878'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.44%  '
